$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<they>"
$ws.Range("C2").Value = 50

# Row 3
$ws.Range("B3").Value = "<on>"

# Row 4
$ws.Range("B4").Value = "<his>"
$ws.Range("C4").Value = 52

# Row 5
$ws.Range("B5").Value = "<there>"
$ws.Range("C5").Value = 53

# Row 6
$ws.Range("B6").Value = "<six>"
$ws.Range("C6").Value = 56

# Row 7
$ws.Range("B7").Value = "<copy>"
$ws.Range("C7").Value = 55

# Row 8
$ws.Range("B8").Value = "<ulo>"
$ws.Range("C8").Value = 53

# Row 9
$ws.Range("B9").Value = "<wonward>"
$ws.Range("C9").Value = 58

# Row 10
$ws.Range("B10").Value = "<delete>"

# Row 11
$ws.Range("B11").Value = "<echo>"
$ws.Range("C11").Value = 51

# Row 12
$ws.Range("B12").Value = "<ta>"

# Row 13
$ws.Range("B13").Value = "<on>"
$ws.Range("C13").Value = 53

# Row 15
$ws.Range("B15").Value = "<first>"
$ws.Range("C15").Value = 55

# Row 16
$ws.Range("C16").Value = 27
